$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.116.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.60%  "
$ws.Range("D3").Value = "'2.452.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("D5").Value = "'309.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Value = "'92.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.11%  "
$ws.Range("E7").Value = "  -2.95%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  -5.73%  "
$ws.Range("D10").Value = "'32.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.82%  "
$ws.Range("D11").Value = "'0.0773"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.86%  "
$ws.Range("E12").Value = "  -0.86%  "
$ws.Range("D13").Value = "'6.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.95%  "
$ws.Range("D14").Value = "'2.832.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.91%  "
$ws.Range("D15").Value = "'2.442.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.38%  "
$ws.Range("D16").Value = "'14.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.11%  "
$ws.Range("D17").Value = "'0.774"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.19%  "
$ws.Range("D18").Value = "'41.078.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.66%  "
$ws.Range("D19").Value = "'6.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.14%  "
$ws.Range("D20").Value = "'0.0₃0909"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.20%  "
$ws.Range("D21").Value = "'11.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.44%  "
$ws.Range("D22").Value = "'67.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.88%  "
$ws.Range("D23").Value = "'234.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.54%  "
$ws.Range("E24").Value = "  -5.07%  "
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("E26").Value = "  -7.03%  "
$ws.Range("D27").Value = "'23.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.10%  "
$ws.Range("E28").Value = "  -5.88%  "
$ws.Range("D29").Value = "'9.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.20%  "
$ws.Range("D30").Value = "'35.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.63%  "
$ws.Range("D31").Value = "'150.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.57%  "
$ws.Range("D32").Value = "'5.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.62%  "
$ws.Range("D33").Value = "'2.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.43%  "
$ws.Range("E34").Value = "  -3.42%  "
$ws.Range("E35").Value = "  -6.40%  "
$ws.Range("D36").Value = "'2.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.74%  "
$ws.Range("D37").Value = "'16.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.20%  "
$ws.Range("E38").Value = "  -6.57%  "
$ws.Range("E39").Value = "  -3.82%  "
$ws.Range("E40").Value = "  -8.34%  "
$ws.Range("D41").Value = "'4.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.78%  "
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").Value = "'19.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.65%  "
$ws.Range("D44").Value = "'1.956.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.36%  "
$ws.Range("E45").Value = "  -6.74%  "
$ws.Range("D46").Value = "'2.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.20%  "
$ws.Range("D47").Value = "'8.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.75%  "
$ws.Range("D48").Value = "'69.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.90%  "
$ws.Range("D49").Value = "'95.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = "  -7.42%  "
$ws.Range("D51").Value = "'73.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.35%  "
